$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "-v" short command moves from the --verbose row (row 2) down to the
# --slicesV row (row 10), and three brand new short commands are added:
#   -h  for --slicesH (row 9)
#   -s  for --searchSubdirs (row 5)
#   -p  for --searchPattern (row 6)
#
# Order matters here: new shared-string entries are appended in the order
# the cells are written, so we add -h, -s, -p (matching the target shared
# string order 43=-h,44=-s,45=-p) before touching the Description text in
# E14 (which becomes the next new shared string, 46) and before moving
# "-v" into its new home (which reuses the existing string, index 4).

# Remove "-v" from its old location (--verbose row).
$ws.Range("A2").ClearContents()

# Add the new short commands. Typing them with a leading apostrophe mimics
# a literal leading "-" being entered as text in Excel, which is what
# produces the quote-prefixed text style seen on these cells.
$ws.Range("A9").Formula = "'-h"
$ws.Range("A5").Formula = "'-s"
$ws.Range("A6").Formula = "'-p"

# Move "-v" to its new row (--slicesV). This reuses the plain (non quote
# prefixed) style already present on the cell.
$ws.Range("A10").Value = "-v"

# Update the --rmExts description so the glob example uses "*.sprite"
# instead of ".sprite".
$ws.Range("E14").Value = "Remove all input files' extensions. Yields files formated as ""*.sprite"""
